$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the new text values. The order below matches the order the
# strings were typed into the workbook (and therefore the order they
# land in xl/sharedStrings.xml): do, some, changes, in, "excel ",
# "as ", well, more.
$ws.Range("C11").Value = "do"
$ws.Range("F6").Value = "some"
$ws.Range("F14").Value = "changes"
$ws.Range("B16").Value = "in"
$ws.Range("M1").Value = "excel "
$ws.Range("L13").Value = "as "
$ws.Range("O8").Value = "well"
$ws.Range("S10").Value = "more"
$ws.Range("P19").Value = "more"

# Highlight S10 and P19 with a yellow fill (matches the new fill/cellXf
# added to styles.xml).
$ws.Range("S10").Interior.Color = 65535
$ws.Range("P19").Interior.Color = 65535

# Leave the selection on S10, like the saved worksheet view.
$ws.Range("S10").Select() | Out-Null
